# Add a new "Conclusion" bullet to the Research section's list, and
# rename the list's former trailing "Conclusion" item to "Bibliography".
#
# Before:  ... Purpose, Method, Results, Conclusion
# After:   ... Purpose, Method, Results, Conclusion, Bibliography

$d = $word.ActiveDocument

# Locate the paragraph that is the final "Conclusion" item (the list
# under the "Research" heading ends the document, so its "Conclusion"
# bullet is the very last paragraph). Using Trim() sidesteps the
# trailing paragraph-mark character that Range.Text always carries, and
# re-scanning the whole body (rather than stopping at the first hit)
# means we land on the *last* paragraph whose text is "Conclusion" -
# there's an earlier, unrelated "Conclusion " bullet (with a trailing
# space) elsewhere in the document that we must not touch.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.Trim() -eq "Conclusion") {
        $targetIndex = $i
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find the trailing 'Conclusion' paragraph"
}

# The paragraph right before it ("Results") ends with its own paragraph
# mark; splitting the text one character before that mark (i.e. right
# after "Results") inserts the new paragraph using "Results"'s own
# (unadorned) paragraph formatting for both halves, instead of dragging
# along the explicit <w:spacing> override that only the final
# "Conclusion"/"Bibliography" paragraph should keep.
$prevParagraph = $d.Paragraphs.Item($targetIndex - 1)
$splitPosition = $prevParagraph.Range.End - 1

$insertRange = $d.Range($splitPosition, $splitPosition)
$insertRange.InsertAfter("`rConclusion")

# The original "Conclusion" paragraph has shifted down by one; rename it.
$lastParagraph = $d.Paragraphs.Item($targetIndex + 1)
$lastParagraph.Range.Text = "Bibliography"
